# Rename isotope labels in the "Name" column (column B) to their proper
# atomic names: "48Ca" -> "48Cd" and "132Sn" -> "132Cs".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Walk every used row in column B (the "Name" column), replacing the
# isotope prefixes wherever they occur at the start of the label. The
# "132Sn" -> "132Cs" rename is applied first (matching the author's
# original edit order), then "48Ca" -> "48Cd".
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row  # xlUp = -4162

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and $val -like "132Sn*") {
        $cell.Value = $val -replace "^132Sn", "132Cs"
    }
}

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and $val -like "48Ca*") {
        $cell.Value = $val -replace "^48Ca", "48Cd"
    }
}

# Reflect the author's final scroll position / cursor selection when the
# workbook was saved.
$win = $excel.ActiveWindow
$win.ScrollRow = 27
$win.ScrollColumn = 1
$null = $ws.Range("O62").Select()
